$d = $word.ActiveDocument

# --- Swahili translation -> English text replacements (table cells / body) ---
$replacements = @(
    @{ Old = "Kichwa cha Video"; New = "Video Title" },
    @{ Old = "Mada"; New = "Topic" },
    @{ Old = "Malengo"; New = "Aim(s)" },
    @{ Old = "Urefu"; New = "Length" },
    @{ Old = "Mahali pa Kambi"; New = "Camp Location" },
    @{ Old = "Wawezeshaji"; New = "Facilitators" },
    @{ Old = "N. ya wanafunzi"; New = "N. of students" },
    @{ Old = "Tarehe"; New = "Date" },
    @{ Old = "Rasilimali"; New = "Resources" },
    @{ Old = "inahitajika"; New = "needed" },
    @{ Old = "Maandalizi"; New = "Preparations" },
    @{ Old = "Muda wa video"; New = "Video time" },
    @{ Old = "Mwezeshaji anafanya nini"; New = "What facilitator does" },
    @{ Old = "Wanachofanya wanafunzi"; New = "What learners do" },
    @{ Old = "Utangulizi Mkuu wa Video ya VMC"; New = "General VMC Video Introduction" },
    @{ Old = "Utangulizi wa Video"; New = "Video Introduction" },
    @{ Old = "Kitendawili"; New = "Riddle" },
    @{ Old = "Kusaidia mchakato, kuchochea mawazo"; New = "Assist the process, provoke thoughts" },
    @{ Old = "Suluhisho"; New = "Solution" }
)

foreach ($pair in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($pair.Old, $true, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null
}

# --- Default document language: Swahili (Kenya) -> Swahili (Tanzania) ---
$d.Styles("Normal").LanguageID = "sw-TZ"
